$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 270; existing rows 270-280 shift down to 271-281,
# matching the weekly update that adds one new Kiwi price record.
$ws.Rows.Item(270).Insert()

# Populate the new row 270 with the new weekly record.
$ws.Cells.Item(270, 1).Value2 = 5
$ws.Cells.Item(270, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(270, 3).Value = "Maule"
$ws.Cells.Item(270, 4).Value2 = 44753
$ws.Cells.Item(270, 5).Value2 = 7
$ws.Cells.Item(270, 6).Value = "Fruta"
$ws.Cells.Item(270, 7).Value2 = 100101
$ws.Cells.Item(270, 8).Value = "Berries"
$ws.Cells.Item(270, 9).Value2 = 100101007
$ws.Cells.Item(270, 10).Value = "Kiwi"
$ws.Cells.Item(270, 11).Value = "Hayward"
$ws.Cells.Item(270, 12).Value = "Primera"
$ws.Cells.Item(270, 13).Value2 = 300
$ws.Cells.Item(270, 14).Value2 = 6000
$ws.Cells.Item(270, 15).Value2 = 6000
$ws.Cells.Item(270, 16).Value2 = 6000
$ws.Cells.Item(270, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(270, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(270, 19).Value2 = 333
$ws.Cells.Item(270, 20).Value2 = 18
